# Apply the text corrections to the country-name column (col B) on Tab28.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab28")

$ws.Range("B4").Value  = "Botswana"               # was "Botswana*"
$ws.Range("B11").Value = "Zambie"                 # was "Zambie*"
$ws.Range("B18").Value = "République du Congo*"   # was "Congo*"
$ws.Range("B19").Value = "RD Congo"                # was "DR Congo*"
$ws.Range("B22").Value = "Sao Tomé-et-Principe"    # was "São Tomé and Príncipe"
$ws.Range("B34").Value = "Soudan du Sud"          # was "Soudan du Sud*"
$ws.Range("B36").Value = "Tanzanie"               # was "UR of Tanzania"
$ws.Range("B54").Value = "Liberia"                # was "Libéria"
$ws.Range("B57").Value = "Nigeria*"               # was "Nigéria*"
$ws.Range("B60").Value = "Togo"                   # was "Togo*"

# Bump the saved workbook window height (bookViews/workbookView@windowHeight).
try {
    $excel.ActiveWindow.Height = 12490
} catch {}
try {
    $wb.Windows.Item(1).Height = 12490
} catch {}
